$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = '[Analysis of myocardial injury in patients with COVID-19 and association between concomitant cardiovascular diseases and severity of COVID-19].'
$ws.Range("D2").Value = 'Objective: To evaluate the cardiovascular damage of patients with COVID-19, and determine the correlation of serum N-terminal pro B-type natriuretic peptide (NT-proBNP) and cardiac troponin-I (cTnI) with the severity of COVID-19, and the impact of concomitant cardiovascular disease on severity of COVID-19 was also evaluated.
 Methods: A cross-sectional study was designed on 150 consecutive patients with COVID-19 in the fever clinic of Tongji Hospital in Wuhan from January to February in 2020, including 126 mild cases and 24 cases in critical care.
 Both univariate and multivariate logistic regression were used to analyze the correlation of past medical history including hypertension, diabetes and coronary heart disease (CHD) , as well as the levels of serum NT-proBNP and cTnI to the disease severity of COVID-19 patients.
 Results: Age, hypersensitive C-reactive protein(hs-CRP) and serum creatinine levels of the patients were higher in critical care cases than in mild cases(all P<0.05).
 Prevalence of male, elevated NT-proBNP and cTnI, hypertension and coronary heart disease were significantly higher in critical cases care patients than in the mild cases(all P<0.05).
 Univariate logistic regression analysis showed that age, male, elevated NT-proBNP, elevated cTnI, elevated hs-CRP, elevated serum creatinine, hypertension, and CHD were significantly correlated with critical disease status(all P<0.05).
 Multivariate logistic regression analysis showed that elevated cTnI(OR=26.909, 95%CI 4.086-177.226, P=0.001) and CHD (OR=16.609, 95%CI 2.288-120.577, P=0.005) were the independent risk factors of critical disease status.
 Conclusions: COVID-19 can significantly affect the heart function and lead to myocardial injury.
 The past medical history of CHD and increased level of cTnI are two independent determinants of clinical disease status in patients with COVID-19.'
$ws.Range("E2").Value = '[ C%Chen%elasticNoEmail%2,   J T% Yan%elasticNoEmail%2,   N% Zhou%elasticNoEmail%2,   J P% Zhao%elasticNoEmail%2,   D W% Wang%elasticNoEmail%2]'
$ws.Range("F2").Value = '10.3760/cma.j.cn112148-20200225-00123'
$ws.Range("G2").Value = 'ELASTIC'
$ws.Range("H2").Value = '2020-03-06'

# Row 3
$ws.Range("C3").Value = '[Diagnosis and treatment of COVID-19: acute kidney injury cannot be ignored].'
$ws.Range("D3").Value = '2019年12月新型冠状病毒肺炎（COVID-19）疫情肆虐武汉，不仅席卷整个中国，还蔓延至全球多个国家和地区，被世界卫生组织（WHO）定义为国际关注的突发公共卫生事件。流行病学资料显示此次COVID-19重症比例高达20%以上。重症COVID-19患者除呼吸系统为病毒的主要靶器官外，肾脏也是主要受累器官之一，且COVID-19合并急性肾损伤（AKI）是患者预后不良的独立危险因素。因此救治COVID-19患者的临床实践中应做到规范化的AKI集束化防治（5R原则）：风险筛查（Risk）、早期识别（Recognition)、及时处理（Response）、肾脏替代治疗（Renal replacement theraphy）及肾脏康复（Recovery）来改善患者的预后。.
'
$ws.Range("E3").Value = '[ X H%Yang%elasticNoEmail%2,   R H% Sun%elasticNoEmail%2,   D C% Chen%elasticNoEmail%2]'
$ws.Range("F3").Value = '10.3760/cma.j.cn112137-20200229-00520'
$ws.Range("G3").Value = 'ELASTIC'
$ws.Range("H3").Value = '2020-03-08'

# Row 5
$ws.Range("E5").Value = '[Chaolin%Huang%NULL%0,                          Yeming%Wang%NULL%0,                          Xingwang%Li%NULL%0,                          Lili%Ren%NULL%0,                          Jianping%Zhao%NULL%0,                          Yi%Hu%NULL%0,                          Li%Zhang%NULL%0,                          Guohui%Fan%NULL%0,                          Jiuyang%Xu%NULL%0,                          Xiaoying%Gu%NULL%0,                          Zhenshun%Cheng%NULL%0,                          Ting%Yu%NULL%0,                          Jiaan%Xia%NULL%0,                          Yuan%Wei%NULL%0,                          Wenjuan%Wu%NULL%0,                          Xuelei%Xie%NULL%0,                          Wen%Yin%NULL%0,                          Hui%Li%NULL%0,                          Min%Liu%NULL%0,                          Yan%Xiao%NULL%0,                          Hong%Gao%NULL%0,                          Li%Guo%NULL%0,                          Jungang%Xie%NULL%0,                          Guangfa%Wang%NULL%0,                          Rongmeng%Jiang%NULL%0,                          Zhancheng%Gao%NULL%0,                          Qi%Jin%NULL%0,                          Jianwei%Wang%wangjw28@163.com%0,                          Bin%Cao%caobin_ben@163.com%0]'

# Row 6
$ws.Range("E6").Value = '[Yan%Deng%NULL%5,                          Wei%Liu%NULL%0,                          Kui%Liu%NULL%2,                          Yuan-Yuan%Fang%NULL%0,                          Jin%Shang%NULL%1,                          Ling%Zhou%NULL%1,                          Ke%Wang%NULL%0,                          Fan%Leng%NULL%1,                          Shuang%Wei%NULL%1,                          Lei%Chen%NULL%1,                          Hui-Guo%Liu%NULL%0,                          Pei-Fang%Wei%NULL%0,                          Pei-Fang%Wei%NULL%0]'

# Row 7
$ws.Range("E7").Value = '[Yichun%Cheng%NULL%1,                          Ran%Luo%NULL%1,                          Kun%Wang%NULL%2,                          Meng%Zhang%NULL%1,                          Zhixiang%Wang%NULL%1,                          Lei%Dong%NULL%1,                          Junhua%Li%NULL%2,                          Ying%Yao%NULL%1,                          Shuwang%Ge%geshuwang@tjh.tjmu.edu.cn%1,                          Gang%Xu%xugang@tjh.tjmu.edu.cn%1]'

# Row 8
$ws.Range("E8").Value = '[Carly%Eastin%NULL%2,                          Travis%Eastin%NULL%1]'

# Row 9
$ws.Range("E9").Value = '[Aguilar%Medina Jos\u00e9 M.%coreGivesNoEmail%5,                        Claro%Vald\u00e9s Ram\u00f3n%coreGivesNoEmail%2,                        Obreg\u00f3n%Santos Angel G.%coreGivesNoEmail%2,                        Proh\u00edas%Mart\u00ednez Juan%coreGivesNoEmail%2,                        Rodr\u00edguez%Blanco Suilbert%coreGivesNoEmail%2]'

# Row 11
$ws.Range("E11").Value = '[Luwen%Wang%NULL%2,                          Xun%Li%NULL%2,                          Hui%Chen%NULL%1,                          Shaonan%Yan%NULL%1,                          Dong%Li%NULL%1,                          Yan%Li%NULL%1,                          Zuojiong%Gong%NULL%1]'

# Row 12
$ws.Range("E12").Value = '[Xiao-Wei%Xu%NULL%0,                          Xiao-Xin%Wu%NULL%0,                          Xian-Gao%Jiang%NULL%0,                          Kai-Jin%Xu%NULL%0,                          Ling-Jun%Ying%NULL%0,                          Chun-Lian%Ma%NULL%0,                          Shi-Bo%Li%NULL%0,                          Hua-Ying%Wang%NULL%0,                          Sheng%Zhang%NULL%0,                          Hai-Nv%Gao%NULL%0,                          Ji-Fang%Sheng%NULL%0,                          Hong-Liu%Cai%NULL%0,                          Yun-Qing%Qiu%NULL%0,                          Lan-Juan%Li%NULL%0]'

# Row 13
$ws.Range("C13").Value = '[Diagnosis and treatment of COVID-19: acute kidney injury cannot be ignored].'
$ws.Range("D13").Value = '2019年12月新型冠状病毒肺炎（COVID-19）疫情肆虐武汉，不仅席卷整个中国，还蔓延至全球多个国家和地区，被世界卫生组织（WHO）定义为国际关注的突发公共卫生事件。流行病学资料显示此次COVID-19重症比例高达20%以上。重症COVID-19患者除呼吸系统为病毒的主要靶器官外，肾脏也是主要受累器官之一，且COVID-19合并急性肾损伤（AKI）是患者预后不良的独立危险因素。因此救治COVID-19患者的临床实践中应做到规范化的AKI集束化防治（5R原则）：风险筛查（Risk）、早期识别（Recognition)、及时处理（Response）、肾脏替代治疗（Renal replacement theraphy）及肾脏康复（Recovery）来改善患者的预后。.
'
$ws.Range("E13").Value = '[ X H%Yang%elasticNoEmail%0,   R H% Sun%elasticNoEmail%0,   D C% Chen%elasticNoEmail%0]'
$ws.Range("F13").Value = '10.3760/cma.j.cn112137-20200229-00520'
$ws.Range("G13").Value = 'ELASTIC'
$ws.Range("H13").Value = '2020-03-08'

# Row 14
$ws.Range("E14").Value = '[Gemin%Zhang%NULL%3,                          Jie%Zhang%945128911@qq.com%1,                          Bowen%Wang%NULL%1,                          Xionglin%Zhu%NULL%1,                          Qiang%Wang%NULL%1,                          Shiming%Qiu%NULL%1]'

# Row 15
$ws.Range("E15").Value = '[Xiaoli%Zhang%NULL%2,                          Huan%Cai%NULL%1,                          Jianhua%Hu%NULL%1,                          Jiangshan%Lian%NULL%1,                          Jueqing%Gu%NULL%1,                          Shanyan%Zhang%NULL%1,                          Chanyuan%Ye%NULL%1,                          Yingfeng%Lu%NULL%1,                          Ciliang%Jin%NULL%1,                          Guodong%Yu%NULL%1,                          Hongyu%Jia%NULL%1,                          Yimin%Zhang%NULL%1,                          Jifang%Sheng%jifang_sheng@zju.edu.cn%1,                          Lanjuan%Li%ljli@zju.edu.cn%0,                          Yida%Yang%yidayang65@zju.edu.cn%1]'

# Row 16
$ws.Range("E16").Value = '[Fei%Zhou%NULL%0,                          Ting%Yu%NULL%0,                          Ronghui%Du%NULL%2,                          Guohui%Fan%NULL%0,                          Ying%Liu%NULL%0,                          Zhibo%Liu%NULL%2,                          Jie%Xiang%NULL%0,                          Yeming%Wang%NULL%0,                          Bin%Song%NULL%2,                          Xiaoying%Gu%NULL%0,                          Lulu%Guan%NULL%2,                          Yuan%Wei%NULL%0,                          Hui%Li%NULL%0,                          Xudong%Wu%NULL%2,                          Jiuyang%Xu%NULL%0,                          Shengjin%Tu%NULL%2,                          Yi%Zhang%NULL%2,                          Hua%Chen%NULL%2,                          Bin%Cao%NULL%0]'

# Row 17
$ws.Range("C17").Value = '[Clinical characteristics and outcomes of 112 cardiovascular disease patients infected by 2019-nCoV].'
$ws.Range("D17").Value = 'Objective: To explore the clinical characteristics and prognosis of the new coronavirus 2019-nCoV patients combined with cardiovascular disease (CVD).
 Methods: A retrospective analysis was performed on 112 COVID-19 patients with CVD admitted to the western district of Union Hospital in Wuhan, from January 20, 2020 to February 15, 2020. They were divided into critical group (ICU, n=16) and general group (n=96) according to the severity of the disease and patients were followed up to the clinical endpoint.
 The observation indicators included total blood count, C-reactive protein (CRP), arterial blood gas analysis, myocardial injury markers, coagulation function, liver and kidney function, electrolyte, procalcitonin (PCT), B-type natriuretic peptide (BNP), blood lipid, pulmonary CT and pathogen detection.
 Results: Compared with the general group, the lymphocyte count (0.74×10(9) (0.34×10(9), 0.94×10(9))/L vs.
 0.99×10(9) (0.71×10(9), 1.29×10(9))/L, P=0.03) was extremely lower in the critical group, CRP (106.98 (81.57, 135.76) mg/L vs.
 34.34 (9.55,76.54) mg/L, P 25 kg/m(2), which was significantly higher than that of survivors (18.95% (18/95), P0.05).
 Conclusion: COVID-19 patients combined with CVD are associated with a higher risk of mortality.
 Critical patients are characterized with lower lymphocyte counts.
 Higher BMI are more often seen in critical patients and non-survivor.
 ACEI/ARB use does not affect the morbidity and mortality of COVID-19 combined with CVD.
 Aggravating causes of death include fulminant inflammation, lactic acid accumulation and thrombotic events.
'
$ws.Range("E17").Value = '[ Y D%Peng%elasticNoEmail%2,   K% Meng%elasticNoEmail%2,   H Q% Guan%elasticNoEmail%2,   L% Leng%elasticNoEmail%2,   R R% Zhu%elasticNoEmail%2,   B Y% Wang%elasticNoEmail%2,   M A% He%elasticNoEmail%2,   L X% Cheng%elasticNoEmail%2,   K% Huang%elasticNoEmail%2,   Q T% Zeng%elasticNoEmail%2]'
$ws.Range("F17").Value = '10.3760/cma.j.cn112148-20200220-00105'
$ws.Range("G17").Value = 'ELASTIC'
$ws.Range("H17").Value = '2020-03-02'
